# fhir ig initial setup
$wb = $excel.ActiveWorkbook

# --- Metadata sheet: bump the generation Date ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Cells.Item(8, 2).Value = "2025-08-20T08:30:34+05:45"

# --- Include #0 sheet: rebuild as a "concept is-a <code>" filter table ---
$inc = $wb.Worksheets.Item("Include #0")

# New 3-column / 4-row layout (was 2-column / 6-row).
$inc.Cells.Item(1, 1).Value = "Property"
$inc.Cells.Item(1, 2).Value = "Operation"

# C1 is a brand-new cell outside the old A1:B6 range, so it has no style of
# its own yet -- borrow the header look from A1/B1 via a formats-only paste.
$inc.Cells.Item(1, 1).Copy()
$inc.Cells.Item(1, 3).PasteSpecial(-4122)
$inc.Cells.Item(1, 3).Value = "Value"

$inc.Cells.Item(2, 1).Value = "concept"
$inc.Cells.Item(2, 2).Value = "is-a"

# C2 is also brand-new -- borrow the body look from B2 the same way. Enter
# the SNOMED code with a leading apostrophe so it is kept as text (not
# coerced to a number), then re-paste the plain body format over it so the
# quote-prefix styling doesn't linger.
$inc.Cells.Item(2, 2).Copy()
$inc.Cells.Item(2, 3).PasteSpecial(-4122)
$inc.Cells.Item(2, 3).Value = "'185389009"
$inc.Cells.Item(2, 2).Copy()
$inc.Cells.Item(2, 3).PasteSpecial(-4122)

$inc.Cells.Item(3, 1).Value = ""
$inc.Cells.Item(3, 2).Value = ""

$inc.Cells.Item(4, 1).Value = "System URI"
$inc.Cells.Item(4, 2).Value = "http://snomed.info/sct"

# The old layout had two extra rows (5 and 6) that no longer exist.
$inc.Range("A5:B6").Clear()

$excel.CutCopyMode = $false
